$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.368.62"
$ws.Range("D3").Value = "1.874.65"
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7149"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.61%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3111"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.23%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07765"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.14"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.61%  "
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("D12").Value = "1.884.96"
$ws.Range("E12").Value = "  +1.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.257"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.53%  "
$ws.Range("E14").Value = "  +0.66%  "
$ws.Range("E15").Value = "  +0.22%  "
$ws.Range("D16").Value = "29.373.96"
$ws.Range("E16").Value = "  +0.56%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.091"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008243"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "240.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.86%  "
$ws.Range("E20").Value = "  +1.28%  "
$ws.Range("D21").Value = "2.124.39"
$ws.Range("E21").Value = "  +0.73%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.791"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.69%  "
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1596"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.87%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.054"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.57%  "
$ws.Range("E29").Value = "  +0.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.423"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.331"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.284"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05312"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.939"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.42%  "
$ws.Range("E35").Value = "  +1.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7412"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.701"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01872"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.63%  "
$ws.Range("D39").Value = "1.231.04"
$ws.Range("E39").Value = "  +5.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.729"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.528"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "110.73"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8891"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.72%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "73.05"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.58%  "
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").Value = "2.021.57"
$ws.Range("E46").Value = "  +0.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.812"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5215"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.80%  "
$ws.Range("E49").Value = "  +2.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.461"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4319"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.40%  "
